$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" section
# (contains "LOB1004"). The footer block that follows it -
# a blank paragraph, the "Ver no Jupiter..." line, and the
# "(c) 2020 ..." copyright line - is being removed from the page.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1004*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $startPara = $anchorIndex + 1
    $endPara = $anchorIndex + 3

    $startPos = $d.Paragraphs.Item($startPara).Range.Start
    $endPos = $d.Paragraphs.Item($endPara).Range.End

    $r = $d.Range($startPos, $endPos)
    $r.Delete()
}
